# Update countries & provincias Spain
# - Reorder "Finlandia" ahead of "Panama" in the country list, and
#   "Togo" ahead of "Mali", updating the Finlandia/Togo stats to their
#   newer figures (Panama/Mali rows keep their previous numbers).
# - Refresh a handful of per-country stat rows (Austria, Albania, Malta)
#   with newer totals.
# - Bump the "Datos actualizados" timestamp cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp header
$ws.Range("A1").Value = "Datos actualizados a 10 de Abril de 2020 a las 12:52"

# Row 19: Austria
$ws.Cells.Item(19, 2).Value = 13398
$ws.Cells.Item(19, 3).Value = 154
$ws.Cells.Item(19, 5).Value = 7015

# Row 45 now becomes "Finlandia" (previously "Panama"), with refreshed figures
$ws.Cells.Item(45, 1).Value = "Finlandia"
$ws.Cells.Item(45, 2).Value = 2769
$ws.Cells.Item(45, 3).Value = 164
$ws.Cells.Item(45, 4).Value = 300
$ws.Cells.Item(45, 5).Value = 2427
$ws.Cells.Item(45, 6).Value = 82
$ws.Cells.Item(45, 7).Value = 0
$ws.Cells.Item(45, 8).Value = 42

# Row 46 now becomes "Panama" (previously "Finlandia"), keeping the old
# Panama figures
$ws.Cells.Item(46, 1).Value = "Panama"
$ws.Cells.Item(46, 2).Value = 2752
$ws.Cells.Item(46, 3).Value = 0
$ws.Cells.Item(46, 4).Value = 16
$ws.Cells.Item(46, 5).Value = 2670
$ws.Cells.Item(46, 6).Value = 107
$ws.Cells.Item(46, 7).Value = 0
$ws.Cells.Item(46, 8).Value = 66

# Row 95: Albania
$ws.Cells.Item(95, 6).Value = 7

# Row 102: Malta
$ws.Cells.Item(102, 2).Value = 350
$ws.Cells.Item(102, 3).Value = 13
$ws.Cells.Item(102, 5).Value = 332

# Row 137 now becomes "Togo" (previously "Mali"), with refreshed figures
$ws.Cells.Item(137, 1).Value = "Togo"
$ws.Cells.Item(137, 2).Value = 76
$ws.Cells.Item(137, 3).Value = 3
$ws.Cells.Item(137, 4).Value = 25
$ws.Cells.Item(137, 5).Value = 48
$ws.Cells.Item(137, 6).Value = 0
$ws.Cells.Item(137, 7).Value = 0
$ws.Cells.Item(137, 8).Value = 3

# Row 138 now becomes "Mali" (previously "Togo"), keeping the old Mali
# figures
$ws.Cells.Item(138, 1).Value = "Mali"
$ws.Cells.Item(138, 2).Value = 74
$ws.Cells.Item(138, 3).Value = 0
$ws.Cells.Item(138, 4).Value = 22
$ws.Cells.Item(138, 5).Value = 45
$ws.Cells.Item(138, 6).Value = 0
$ws.Cells.Item(138, 7).Value = 0
$ws.Cells.Item(138, 8).Value = 7
